$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# Insert 4 new data rows, in top-down order, using the existing row 7
# as the format template (PasteSpecial formats only keeps the style
# indices already present in the sheet instead of minting new ones).
# -----------------------------------------------------------------

# 1) DECLOPHEN -> goes right after row 7 (BI-PROFENID)
$ws.Rows("8:8").Insert()
$ws.Range("A7:Q7").Copy()
$ws.Range("A8:Q8").PasteSpecial(-4122)

# 2) FUSI -> goes right after FELDENE (now row 10), before MUCO
$ws.Rows("11:11").Insert()
$ws.Range("A7:Q7").Copy()
$ws.Range("A11:Q11").PasteSpecial(-4122)

# 3) VOLTAREN -> goes right after PHYTO K (now row 14), before the total row
$ws.Rows("15:15").Insert()
$ws.Range("A7:Q7").Copy()
$ws.Range("A15:Q15").PasteSpecial(-4122)

# 4) سرنجات 3 سم -> goes right after VOLTAREN, before the total row
$ws.Rows("16:16").Insert()
$ws.Range("A7:Q7").Copy()
$ws.Range("A16:Q16").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# -----------------------------------------------------------------
# Fill in the values for the whole (now 10-row) data block.
# -----------------------------------------------------------------

$rows = @(
    @{ r = 7;  n = 1;  name = "BI-PROFENID 150MG 20 SCORED TABS."; bal = "2:0";  req = "1"; price = "54.00"; sale = "27.0000";  txn = "0:1" },
    @{ r = 8;  n = 2;  name = "DECLOPHEN 75MG/3ML 3 AMPOULES";      bal = "7:2";  req = "1"; price = "36.00"; sale = "11.8800";  txn = "0:1" },
    @{ r = 9;  n = 3;  name = "FAROVIGA 100MG 12 F.C.TAB.";         bal = "1:11"; req = "1"; price = "108.00"; sale = "35.6400"; txn = "0:4" },
    @{ r = 10; n = 4;  name = "FELDENE 20MG/ML I.M. 6 AMP.";        bal = "2:5";  req = "1"; price = "63.00"; sale = "63.0000";  txn = "1:0" },
    @{ r = 11; n = 5;  name = "FUSI 2% CREAM 15 GM";                bal = "0:0";  req = "1"; price = "35.00"; sale = "35.0000";  txn = "1:0" },
    @{ r = 12; n = 6;  name = "MUCO 15MG/5ML SYRUP 100ML";          bal = "1:0";  req = "1"; price = "35.00"; sale = "35.0000";  txn = "1:0" },
    @{ r = 13; n = 7;  name = "ORS 10 SACHET";                      bal = "6:9";  req = "1"; price = "40.00"; sale = "4.0000";   txn = "0:1" },
    @{ r = 14; n = 8;  name = "PHYTO K 10 MG 50 F.C.TAB.";          bal = "1:0";  req = "1"; price = "72.50"; sale = "-14.5000"; txn = "0:-1" },
    @{ r = 15; n = 9;  name = "VOLTAREN 75MG/3ML 3 AMP.";           bal = "1:0";  req = "1"; price = "51.00"; sale = "33.6600";  txn = "0:2" },
    @{ r = 16; n = 10; name = "سرنجات 3 سم";                        bal = "0:0";  req = "0"; price = "2.00";  sale = "6.0000";   txn = "3:0" }
)

foreach ($row in $rows) {
    $r = $row.r
    $ws.Cells.Item($r, 1).Value = $row.n        # A - م (sequence number)
    $ws.Cells.Item($r, 3).Value = $row.name     # C - الاسم
    $ws.Cells.Item($r, 8).Value = $row.bal      # H - الرصيد الحالي
    $ws.Cells.Item($r, 12).Value = $row.req     # L - حد الطلب
    $ws.Cells.Item($r, 14).Value = $row.price   # N - السعر
    $ws.Cells.Item($r, 16).Value = $row.sale    # P - سعر البيع
    $ws.Cells.Item($r, 17).Value = $row.txn     # Q - عدد التعاملات
}

# -----------------------------------------------------------------
# Update the totals row (now row 17) and the footer timestamp (row 18).
# -----------------------------------------------------------------
$ws.Cells.Item(17, 16).Value = 236.68
$ws.Cells.Item(18, 1).Value = "Sunday, 15 June, 2025 11:01 AM"
